$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40"; the edit replaces that value
# with the text "1". Assigning a plain numeric-looking string (e.g.
# Value = "1") would be auto-coerced to the number 1, so instead we
# build it as a text-formula result and flatten it back to a literal
# value in place - this keeps the cell's existing formatting/style
# untouched while still storing "1" as text rather than a number.
$cell = $ws.Cells.Item(11, 2)
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
